$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.134.35'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '3.779.82'
$ws.Range('E3').Value = '  +6.48%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').Value = '''137.77'
$ws.Range('E6').Value = '  +4.27%  '
$ws.Range('D7').Value = '3.766.88'
$ws.Range('E7').Value = '  +6.42%  '
$ws.Range('D8').Value = '''0.638'
$ws.Range('E8').Value = '  -3.68%  '
$ws.Range('D9').Value = '''0.998'
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = '''0.758'
$ws.Range('E10').Value = '  -3.24%  '
$ws.Range('E11').Value = '  +5.26%  '
$ws.Range('E12').Value = '  +28.75%  '
$ws.Range('D13').Value = '''42.36'
$ws.Range('E13').Value = '  -2.46%  '
$ws.Range('D14').Value = '''10.23'
$ws.Range('D15').Value = '4.378.79'
$ws.Range('E15').Value = '  +6.18%  '
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').Value = '3.805.95'
$ws.Range('E17').Value = '  +6.75%  '
$ws.Range('D18').Value = '''20.44'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').Value = '''13.36'
$ws.Range('E19').Value = '  +3.76%  '
$ws.Range('D20').Value = '''1.11'
$ws.Range('E20').Value = '  -0.84%  '
$ws.Range('D21').Value = '67.252.19'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('D22').Value = '''436.15'
$ws.Range('E22').Value = '  -3.10%  '
$ws.Range('D23').Value = '''15.08'
$ws.Range('E23').Value = '  +13.95%  '
$ws.Range('D24').Value = '''88.90'
$ws.Range('E24').Value = '  -1.69%  '
$ws.Range('E25').Value = '  -5.99%  '
$ws.Range('D26').Value = '''37.41'
$ws.Range('E26').Value = '  +9.11%  '
$ws.Range('D27').Value = '''3.27'
$ws.Range('E27').Value = '  -3.61%  '
$ws.Range('D28').Value = '''9.73'
$ws.Range('E28').Value = '  -3.07%  '
$ws.Range('D29').Value = '''5.14'
$ws.Range('E29').Value = '  +6.52%  '
$ws.Range('D30').Value = '''12.50'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = '''0.121'
$ws.Range('E31').Value = '  +3.27%  '
$ws.Range('D32').Value = '''2.73'
$ws.Range('E32').Value = '  -2.31%  '
$ws.Range('D33').Value = '''7.14'
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '''0.161'
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = '''41.13'
$ws.Range('E35').Value = '  +4.73%  '
$ws.Range('D36').Value = '''57.75'
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('D37').Value = '''0.998'
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').Value = '''0.0481'
$ws.Range('E38').Value = '  -5.36%  '
$ws.Range('D39').Value = '''2.99'
$ws.Range('E39').Value = '  +27.98%  '
$ws.Range('D40').Value = '''0.144'
$ws.Range('E40').Value = '  -3.64%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0675'
$ws.Range('E41').Value = '  -8.15%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '''0.999'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = '''3.37'
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '''146.97'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''26.49'
$ws.Range('E45').Value = '  +20.85%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '''2.08'
$ws.Range('E46').Value = '  +3.54%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '''4.32'
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '''2.84'
$ws.Range('E48').Value = '  -6.40%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '''3.06'
$ws.Range('E49').Value = '  +19.03%  '
$ws.Range('D50').Value = '''2.58'
$ws.Range('E50').Value = '  -7.43%  '
$ws.Range('D51').Value = '''0.298'
$ws.Range('E51').Value = '  -4.22%  '
